$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6248894929885864
$ws.Range("B1").Value = 0.759450376033783
$ws.Range("C1").Value = 1.04095983505249
$ws.Range("D1").Value = 3.27154016494751
$ws.Range("E1").Value = 6.345268249511719
